# Duplicate slide 14 (the "Database / Application" figure) so the copy
# becomes the new slide 15, then update the duplicate's "Query Engine"
# label to "Inference and Query Engine". This pushes the former slides
# 15 and 16 down to become slides 16 and 17.

$p = $ppt.ActivePresentation

$source = $p.Slides.Item(14)
$dup = $source.Duplicate()

# Duplicate() returns a SlideRange containing the new slide, which is
# inserted immediately after the source slide (i.e. at index 15).
$newSlide = $dup.Item(1)

foreach ($shape in $newSlide.Shapes) {
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Query Engine") {
            $shape.TextFrame.TextRange.Text = "Inference and Query Engine"
        }
    }
}
